$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 471, shifting existing rows 471..496 down to 472..497.
$ws.Rows.Item(471).Insert()

# Populate the newly inserted row 471 with the new record.
$ws.Cells.Item(471, 1).Value = 4
$ws.Cells.Item(471, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(471, 3).Value = "Los Lagos"
$ws.Cells.Item(471, 4).Value = 45267
$ws.Cells.Item(471, 5).Value = 10
$ws.Cells.Item(471, 6).Value = 100112021
$ws.Cells.Item(471, 7).Value = "Ají"
$ws.Cells.Item(471, 8).Value = "Inferno"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 120
$ws.Cells.Item(471, 11).Value = 41000
$ws.Cells.Item(471, 12).Value = 41000
$ws.Cells.Item(471, 13).Value = 41000
$ws.Cells.Item(471, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(471, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(471, 16).Value = 4100
$ws.Cells.Item(471, 17).Value = 10
$ws.Cells.Item(471, 18).Value = "Hortaliza"
